# Menu update via menuPublish script
$wb = $excel.ActiveWorkbook

# --- Grab the two worksheets we need by their ORIGINAL names first ---
$wsCarne  = $wb.Worksheets.Item("CARNE & INSALATE")
$wsBurger = $wb.Worksheets.Item("BURGER & FRITTI")

# --- Rename sheets ---
$wsCarne.Name  = "PIATTI UNICI E INSALATE"
$wsBurger.Name = "HAMBURGER E SNACK"

# --- Update text on the "PIATTI UNICI E INSALATE" sheet (was CARNE & INSALATE) ---
$wsCarne.Range("A3").Value = "PIATTI UNICI E INSALATE"
$wsCarne.Range("C4").Value = "MENU OKTOBERFEST"
$wsCarne.Range("G5").Value = "Tre tipologie di wurstel bavaresi: bockwurst, bratwurst e weisswurst, cotti alla griglia e accompagnati con crauti e bretzel, con birra media inclusa."
$wsCarne.Range("G6").Value = "Stinco di maiale servito con salsa a base di birra fatta in casa, accompagnato con crauti e patate al forno, con birra media inclusa."
$wsCarne.Range("G7").Value = "Combinazione dei nostri piatti bavaresi, iniziando dal piatto wurstel per passare allo strinco, con  litro di birra incluso."
$wsCarne.Range("C9").Value = "PIATTI UNICI"

# Row heights for the (now taller) MENU WURSTEL / MENU STINCO rows
$wsCarne.Rows.Item(5).RowHeight = 46.25
$wsCarne.Rows.Item(6).RowHeight = 46.25

# Clear the stray "1" that had been left in Y8
$wsCarne.Range("Y8").ClearContents()

# --- Update text on the "HAMBURGER E SNACK" sheet (was BURGER & FRITTI) ---
$wsBurger.Range("A3").Value = "HAMBURGER E SNACK"

# --- Selections / scroll position ---
# Set selection on the (non-active) HAMBURGER E SNACK sheet first, scrolled back to the top.
$wsBurger.Range("A3").Select() | Out-Null

# Re-activate the PIATTI UNICI E INSALATE sheet (it remains the selected tab)
# and move its selection from G7 to K5.
$wsCarne.Activate() | Out-Null
$wsCarne.Range("K5").Select() | Out-Null

Write-Host "menuPublish edits applied"
